$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.444.88"
$ws.Range("E2").Value = "  -0.49%  "
$ws.Range("D3").Value = "1.864.51"
$ws.Range("E3").Value = "  -0.86%  "
$ws.Range("E4").Value = "  -1.83%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.01"
$ws.Range("E5").Value = "  -1.28%  "
$ws.Range("E6").Value = "  -1.64%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5061"
$ws.Range("E7").Value = "  -1.96%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3908"
$ws.Range("E8").Value = "  -2.08%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08314"
$ws.Range("E9").Value = "  -1.14%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.47"
$ws.Range("E10").Value = "  +0.48%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.101"
$ws.Range("E11").Value = "  -1.49%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.181"
$ws.Range("E12").Value = "  -1.71%  "
$ws.Range("D13").Value = "1.864.31"
$ws.Range("E13").Value = "  +1.39%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.27"
$ws.Range("E14").Value = "  -1.73%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.224"
$ws.Range("E15").Value = "  -0.60%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.009"
$ws.Range("E16").Value = "  -1.99%  "
$ws.Range("E17").Value = "  -1.49%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "91.12"
$ws.Range("E18").Value = "  -0.34%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06721"
$ws.Range("E19").Value = "  -1.10%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.60"
$ws.Range("E20").Value = "  -1.23%  "
$ws.Range("E21").Value = "  -1.70%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.897"
$ws.Range("E22").Value = "  -1.63%  "
$ws.Range("D23").Value = "28.498.08"
$ws.Range("E23").Value = "  -0.27%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.05"
$ws.Range("E24").Value = "  -1.40%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.203"
$ws.Range("E25").Value = "  -4.20%  "
$ws.Range("D26").Value = "2.079.02"
$ws.Range("E26").Value = "  +1.47%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "157.79"
$ws.Range("E27").Value = "  -3.06%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.52"
$ws.Range("E28").Value = "  -1.99%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.407"
$ws.Range("E29").Value = "  +1.08%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "126.26"
$ws.Range("E30").Value = "  -1.42%  "
$ws.Range("E31").Value = "  -1.81%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.032"
$ws.Range("E32").Value = "  -1.14%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.734"
$ws.Range("E33").Value = "  -1.89%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.611"
$ws.Range("E34").Value = "  -1.49%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02439"
$ws.Range("E35").Value = "  -0.29%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06561"
$ws.Range("E36").Value = "  +0.53%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "8.951"
$ws.Range("E37").Value = "  -0.29%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2151"
$ws.Range("E38").Value = "  -2.09%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.004"
$ws.Range("E39").Value = "  -1.32%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.176"
$ws.Range("E40").Value = "  -1.47%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.229"
$ws.Range("E41").Value = "  -4.14%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6332"
$ws.Range("E42").Value = "  -2.12%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.06"
$ws.Range("E43").Value = "  -2.32%  "
$ws.Range("E44").Value = "  -1.45%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5976"
$ws.Range("E45").Value = "  -1.44%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.95"
$ws.Range("E46").Value = "  -1.24%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.680"
$ws.Range("E47").Value = "  -1.74%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.989"
$ws.Range("E48").Value = "  -0.81%  "
$ws.Range("E49").Value = "  -0.42%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.204"
$ws.Range("E50").Value = "  -0.87%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.133"
$ws.Range("E51").Value = "  -7.90%  "
